$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.350.87"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "3.912.58"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "486.21"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "145.77"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "0.0000345"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("D12").Value = "43.12"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "10.75"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").Value = "4.540.50"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "3.907.51"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "14.35"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "20.02"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "68.403.14"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").Value = "432.96"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "15.23"
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").Value = "88.32"
$ws.Range("D25").Value = "11.47"
$ws.Range("E25").Value = "  +15.77%  "
$ws.Range("D26").Value = "11.21"
$ws.Range("E26").Value = "  +10.62%  "
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "37.93"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "5.73"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "720.05"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "13.80"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "2.94"
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("D34").Value = "6.19"
$ws.Range("E34").Value = "  +14.83%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "41.43"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0880"
$ws.Range("E36").Value = "  +3.88%  "
$ws.Range("D37").Value = "60.93"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0498"
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.394"
$ws.Range("E41").Value = "  +16.12%  "
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  +18.14%  "
$ws.Range("D43").Value = "3.11"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  +5.54%  "
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +5.08%  "
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  -5.08%  "
$ws.Range("D50").Value = "145.16"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "0.0₆0342"
$ws.Range("E51").Value = "  +27.58%  "
